# This script applies the weekly data update described by the commit:
# "Fruta / hortaliza, semanal"
#
# Two brand-new observation rows are inserted right before the current
# row 837 (pushing the existing rows 837-897 down to 839-899, and growing
# the sheet's used range from A1:R897 to A1:R899). The two freshly
# inserted rows (837 and 838) are then populated with new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 837 (this shifts rows 837:897 down to 839:899)
$ws.Rows.Item(837).Insert()
$ws.Rows.Item(837).Insert()

# --- New row 837 ---
$ws.Cells.Item(837, 1).Value = 10
$ws.Cells.Item(837, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(837, 3).Value = "La Araucanía"
$ws.Cells.Item(837, 4).Value = 44461
$ws.Cells.Item(837, 5).Value = 9
$ws.Cells.Item(837, 6).Value = 100112020
$ws.Cells.Item(837, 7).Value = "Tomate"
$ws.Cells.Item(837, 8).Value = "Larga vida"
$ws.Cells.Item(837, 9).Value = "Primera"
$ws.Cells.Item(837, 10).Value = 300
$ws.Cells.Item(837, 11).Value = 9000
$ws.Cells.Item(837, 12).Value = 9000
$ws.Cells.Item(837, 13).Value = 9000
$ws.Cells.Item(837, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(837, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(837, 16).Value = 900
$ws.Cells.Item(837, 17).Value = 10
$ws.Cells.Item(837, 18).Value = "Hortaliza"

# --- New row 838 ---
$ws.Cells.Item(838, 1).Value = 10
$ws.Cells.Item(838, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(838, 3).Value = "La Araucanía"
$ws.Cells.Item(838, 4).Value = 44461
$ws.Cells.Item(838, 5).Value = 9
$ws.Cells.Item(838, 6).Value = 100112020
$ws.Cells.Item(838, 7).Value = "Tomate"
$ws.Cells.Item(838, 8).Value = "Larga vida"
$ws.Cells.Item(838, 9).Value = "Segunda"
$ws.Cells.Item(838, 10).Value = 40
$ws.Cells.Item(838, 11).Value = 13000
$ws.Cells.Item(838, 12).Value = 14000
$ws.Cells.Item(838, 13).Value = 13500
$ws.Cells.Item(838, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(838, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(838, 16).Value = 750
$ws.Cells.Item(838, 17).Value = 18
$ws.Cells.Item(838, 18).Value = "Hortaliza"
